$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.575.83'
$ws.Range("E2").Value = '  +0.11%  '

$ws.Range("D3").Value = '1.960.08'
$ws.Range("E3").Value = '  +2.20%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9983'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.99%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9985'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4850'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.81%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2919'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.86%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06745'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.40%  '

$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.30'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.50%  '

$ws.Range("B11").Value = 'Litecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '109.08'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.21%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.963.10'
$ws.Range("E12").Value = '  +2.61%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07762'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.65%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.417'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.74%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6815'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.56%  '

$ws.Range("B16").Value = 'BitcoinCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '290.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.85%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '30.588.81'
$ws.Range("E17").Value = '  +0.19%  '

$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.23%  '

$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").Value = '2.212.75'
$ws.Range("E19").Value = '  +2.31%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007665'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.14%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.619'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.77%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9978'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.22%  '

$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9970'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.36%  '

$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.588'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.25%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.853'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.95%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.05%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.33%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.171'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.02%  '

$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1063'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.435'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.42%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.658'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +14.50%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.440'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.94%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05073'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.11%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7648'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.40%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.176'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.11%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.727'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.38%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02033'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.57%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.710'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.12%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.143'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.354'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.56%  '

$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '109.29'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4444'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.00%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8760'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.49%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '69.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.37%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9981'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.497'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.39%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1272'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.71%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.334'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.16%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.79'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.51%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '927.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.95%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '47.05'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.94%  '
